$s11 = "Invoke browser`nLoad {testUrl}`nType admin in UserNameField with id = user_login`nClear UserNameField with id = user_login`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nGet tagname of LogInBtn with id = wp-submit Assign {LogInBtnTagName}`n{LogInBtnTagName} VerifyEqual input`nGet value of attribute: value of LogInBtn with id = wp-submit  Assign {LogInBtnValue}`n{LogInBtnValue} VerifyEqual Log In`nIs LogInBtn with id = wp-submit displayed Assign {LogInBtnDisplayed}`n{LogInBtnDisplayed} VerifyEqual true`nIs LogInBtn with id = wp-submit enabled Assign {LogInBtnEnabled}`n{LogInBtnEnabled} VerifyEqual true`nGet text of UserNameLabel with xpath = //form[@id='loginform']/p[1]/label Assign {UserNameLabelText}`n{UserNameLabelText} VerifyEqual  Username`nIs RememberMeChkBox with id = rememberme selected Assign {RememberMeSelected}`n{RememberMeSelected} VerifyEqual false`nClick RememberMeChkBox with id = rememberme`nIs RememberMeChkBox with id = rememberme selected Assign {RememberMeSelected}`n{RememberMeSelected} VerifyEqual true`nClick LogInBtn with id = wp-submit`nGet page title Assign {title}`n{title} VerifyEqual Dashboard ‹ test — WordPress`nGet current page url Assign {pageUrl}`n{pageUrl} VerifyEqual http://127.0.0.1/wordpress/wp-admin/`nQuit browser"
$s12 = "Invoke browser`nLoad {testUrl}`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nClick LogInBtn with id = wp-submit`nType test post title in PostTitleField with id = title`nType test post content in PostContentField with id = content`nClick PublishBtn with id = publish`nSleep 2`nIs ViewPostLink with linkText = View post displayed Assign {ViewPostLinkDisplayed}`n{ViewPostLinkDisplayed} VerifyEqual true`nIs EditPostLink with linkText = Edit post displayed Assign {EditPostLinkDisplayed}`n{EditPostLinkDisplayed} VerifyEqual true`nQuit browser"
$s13 = "Invoke browser`nLoad {testUrl}`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nClick LogInBtn with id = wp-submit`nClick PostsLink with xpath = //li[@id='menu-posts']/a/div[3]`nGet page title Assign {title}`n{title} VerifyEqual Posts ‹ test — WordPress`nSelect Edit from ActionsListBox with xpath = //form[@id='posts-filter']/div[1]/div[1]/select`nClick AddedPostLink with xpath = //tr[contains(@id,'post')]/td[1]/strong/a`nGet page title Assign {title}`n{title} VerifyEqual Edit Post ‹ test — WordPress`nClick MoveToTrashLink with xpath = //div[@id='delete-action']/a`nQuit browser`n"
$s15 = "Invoke browser`nLoad {testUrl}`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nClick LogInBtn with id = wp-submit`nIs SettingsSection with id = adv-settings displayed Assign {SettingsSectionDisplayed}`nPrint {SettingsSectionDisplayed}`nIf (not({SettingsSectionDisplayed}))`n  {SettingsSectionDisplayed} VerifyEqual false`n  Click ScreenOptionsBtn with id = show-settings-link`n  Sleep 2`nEnd If`nIs SettingsSection with id = adv-settings displayed Assign {SettingsSectionDisplayed}`nIf ({SettingsSectionDisplayed})`n  {SettingsSectionDisplayed} VerifyEqual true  `nEnd If`nElse`n   Print {SettingsSectionDisplayed}`nEnd Else`nQuit browser"
$s16 = "Invoke browser`nLoad {testUrl}`nType admin in UserNameField with id = user_login`nType test in PasswordField with id = user_pass`nClick LogInBtn with id = wp-submit`nClick ScreenOptionsBtn with id = show-settings-link`nIs SettingsSection with id = adv-settings displayed Assign {SettingsSectionDisplayed}`n#Print {SettingsSectionDisplayed}`nWhile (({SettingsSectionDisplayed}))`n  {SettingsSectionDisplayed} VerifyEqual true  `n  Click ScreenOptionsBtn with id = show-settings-link`n  Is SettingsSection with id = adv-settings displayed Assign {SettingsSectionDisplayed}`nEnd While`n{SettingsSectionDisplayed} VerifyEqual false `nQuit browser"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: update existing rows' script column (testUrl -> {testUrl}) ---
$ws1.Range("C1").Value = $s11
$ws1.Range("C2").Value = $s12
$ws1.Range("C3").Value = $s13

# --- Sheet1 row 5 (new): Enabled TC_005 with new While-loop script ---
# (Set up row 5 and its "TC_005" label before rewriting row 4's script text, so
#  new shared-string entries land in the same order the source workbook uses:
#  TC_004, ...unchanged..., TC_005, If/While script, While-loop script.)
$ws1.Range("A4:C4").Copy()
$ws1.Range("A5:C5").PasteSpecial(-4122)
$ws1.Range("A5").Value = "Enabled"
$ws1.Range("B5").Value = "TC_005"
$ws1.Rows.Item(5).RowHeight = 225

# --- Sheet1 row 4: becomes Disabled TC_004 with updated If/While script, grows taller ---
$ws1.Range("A4").Value = "Disabled"
$ws1.Range("C4").Value = $s15
$ws1.Rows.Item(4).RowHeight = 300

$ws1.Range("C5").Value = $s16

# --- Selections: Sheet2 first (not left active), then Sheet1 last (stays the active tab) ---
[void]$ws2.Range("C1").Select()
[void]$ws1.Range("C1").Select()
